$p = $ppt.ActivePresentation

# --- Slide 17 ("9. Anexos"), body placeholder (shape 2): HU16 bullet text fix ---
$slide = $p.Slides.Item(17)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# The HU16 bullet is the 4th paragraph of this text box:
#   "HU16 (Alta): Como docente evaluador, quiero acceder al uso del sistema para entender su uso."
# The author reselected "acceder al uso del sistema " and retyped it as
# "acceder al sistema ", so only that inner span of characters changes.
$para = $tr.Paragraphs(4, 1)

$oldSpan = "acceder al uso del sistema "
$newSpan = "acceder al sistema "
$startPos = $para.Text.IndexOf($oldSpan) + 1

$span = $para.Characters($startPos, $oldSpan.Length)
$span.Text = $newSpan
